# Generate Report for Handoff
# Updates the localization-status workbook so that the two files that were
# previously "Handed back: in sync with en-US" (acfabba4-... and
# c40d2e80-...) are now shown as "Ready for handoff", refreshes their
# handoff timestamps, and records a "stale handback" Error Detail message
# for each of them on the zh-cn / de-de sheets. Also widens the
# "Error Detail" column so the new message is readable.

$wb = $excel.ActiveWorkbook

$statusText = "Ready for handoff"
$newHandoffDate = "2016-09-04 04:29:45"

$acfErrZhDe = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea1c38b73b40830c63efc51c3cb143365546c2a9/e2e/acfabba4-400e-4ff2-ba75-48ddd1e0bb13.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60577b01da557af59c2c54f76f517b320b9055f6/e2e/acfabba4-400e-4ff2-ba75-48ddd1e0bb13.md."
$c40ErrZhDe = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ea1c38b73b40830c63efc51c3cb143365546c2a9/e2e/c40d2e80-a511-4666-899b-6a01b98736c4.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/60577b01da557af59c2c54f76f517b320b9055f6/e2e/c40d2e80-a511-4666-899b-6a01b98736c4.md."

# ---------------------------------------------------------------------
# Overview sheet: rows 4 (acfabba4-...) and 5 (c40d2e80-...)
#   zh-cn / de-de status columns -> "Ready for handoff"
#   Latest HO Xliff Generate Date -> 2016-09-04 04:29:45
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E4").Value = $statusText
$wsOverview.Range("F4").Value = $statusText
$wsOverview.Range("G4").Value = $newHandoffDate

$wsOverview.Range("E5").Value = $statusText
$wsOverview.Range("F5").Value = $statusText
$wsOverview.Range("G5").Value = $newHandoffDate

# ---------------------------------------------------------------------
# zh-cn sheet: rows 4 (acfabba4-...) and 5 (c40d2e80-...)
#   Status -> "Ready for handoff"
#   Latest Handoff Datetime -> 2016-09-04 04:29:41
#   Error Detail -> stale-handback message
#   Error Detail column (P) width -> 40
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C4").Value = $statusText
$wsZhCn.Range("H4").Value = "2016-09-04 04:29:41"
$wsZhCn.Range("P4").Value = $acfErrZhDe

$wsZhCn.Range("C5").Value = $statusText
$wsZhCn.Range("H5").Value = "2016-09-04 04:29:41"
$wsZhCn.Range("P5").Value = $c40ErrZhDe

$wsZhCn.Columns.Item(16).ColumnWidth = 39.1

# ---------------------------------------------------------------------
# de-de sheet: rows 4 (acfabba4-...) and 5 (c40d2e80-...)
#   Status -> "Ready for handoff"
#   Latest Handoff Datetime -> 2016-09-04 04:29:45
#   Error Detail -> stale-handback message
#   Error Detail column (P) width -> 40
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C4").Value = $statusText
$wsDeDe.Range("H4").Value = $newHandoffDate
$wsDeDe.Range("P4").Value = $acfErrZhDe

$wsDeDe.Range("C5").Value = $statusText
$wsDeDe.Range("H5").Value = $newHandoffDate
$wsDeDe.Range("P5").Value = $c40ErrZhDe

$wsDeDe.Columns.Item(16).ColumnWidth = 39.1
